$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet
$ws.Name = "Gamma2F"

# Correct tiny floating point rounding on C13 / F13
$ws.Range("C13").Value = 0.9963614538049353
$ws.Range("F13").Value = 0.9963614538049353

# Append a new data row (row 16) with the Gaussian Quadrature averages
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(16, 3).Value = 1.363972170137047
$ws.Cells.Item(16, 4).Value = 0.3837816257825609
$ws.Cells.Item(16, 5).Value = 1.050433605735367
$ws.Cells.Item(16, 6).Value = 1.363972170137047
$ws.Cells.Item(16, 7).Value = 0.6828040584107109
$ws.Cells.Item(16, 8).Value = 1.133475496403593
$ws.Cells.Item(16, 9).Value = 1.133176527350012
$ws.Cells.Item(16, 10).Value = 0.3837816257825609
$ws.Cells.Item(16, 11).Value = 0.717107615758964
$ws.Cells.Item(16, 12).Value = 1.040539892948006
$ws.Cells.Item(16, 13).Value = 0.9579405806365485

# Match the bold/bordered/centered style used by column A's header cells
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
